# Adds a new "Dynamic_Factor" row (k_c dynamic factor) into both the
# "summary" and "results" sheets, just above the existing "SC Configuration"
# block, and refreshes the handful of downstream Contact-Force-Spectrum
# results that change because of it.

$wb = $excel.ActiveWorkbook

function Insert-DynamicFactorRow {
    param(
        [object]$ws,
        [int]$row          # row index the new "Dynamic_Factor" row should occupy
    )

    # Push everything at/after $row down by one, carrying values + formats.
    $ws.Rows.Item($row).Insert()

    # The row that is now blank at $row should look like the row below it
    # (which is the original "SC Configuration" row, now shifted down one).
    $srcRow = $ws.Range("A" + ($row + 1) + ":F" + ($row + 1))
    $dstRow = $ws.Range("A" + $row + ":F" + $row)
    $srcRow.Copy()
    $dstRow.PasteSpecial(-4122)  # xlPasteFormats

    # Column B on the new row keeps the un-labelled style from the row above
    # (the merged "Input Variables" block continues through it) rather than
    # the "SC Configuration" label style it just inherited.
    $ws.Range("B" + ($row - 1)).Copy()
    $ws.Range("B" + $row).PasteSpecial(-4122)  # xlPasteFormats

    # Fill in the actual Dynamic_Factor data.
    $ws.Range("C" + $row).Value = "Dynamic_Factor"
    $ws.Range("D" + $row).Value = 1.1
    $ws.Range("E" + $row).Value = 1
    $ws.Range("F" + $row).Value = 1
}

# ---- Sheet "summary": insert before old row 24 -----------------------
$ws1 = $wb.Worksheets.Item("summary")
Insert-DynamicFactorRow $ws1 24

$ws1.Range("B10:B23").UnMerge()
$ws1.Range("B10:B24").Merge()
# "B24:B43" / "A10:A43" were auto-extended by the row insert to
# "B25:B44" / "A10:A44" already.

# ---- Sheet "results": insert before old row 57 ------------------------
$ws2 = $wb.Worksheets.Item("results")
Insert-DynamicFactorRow $ws2 57

$ws2.Range("B43:B56").UnMerge()
$ws2.Range("B43:B57").Merge()
# "B57:B76" / "A43:A76" were auto-extended by the row insert to
# "B58:B77" / "A43:A77" already.

# ---- Downstream Contact-Force-Spectrum-k_c results that shifted with it
$ws2.Range("D7").Value = 0.4590328693389893
$ws2.Range("D8").Value = 54.82754158933363
$ws2.Range("D9").Value = 24.19958007275333

$ws2.Range("D11").Value = 0.491530641913414
$ws2.Range("D12").Value = 58.70912196494103
$ws2.Range("D13").Value = 23.70804736398622

$ws2.Range("D21").Value = 0.1827430725097657
$ws2.Range("D22").Value = 21.82709360795825
$ws2.Range("D23").Value = 31.90132257067649

$ws2.Range("D25").Value = 0.2099792122840881
$ws2.Range("D26").Value = 25.08021704628614
$ws2.Range("D27").Value = 30.59905925286444

$ws2.Range("D35").Value = 0.2177794933319092
$ws2.Range("D36").Value = 1.633346199989319
$ws2.Range("D37").Value = 55.35386401451498

$ws2.Range("D39").Value = 0.3037856727838517
$ws2.Range("D40").Value = 2.278392545878888
$ws2.Range("D41").Value = 50.09367007459512
